$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The portfolio rows for policy 1 (row 2) and policy 2 (row 3) are
# swapped: row 2 now holds policy id 2's data, row 3 holds policy id 1's
# data. Column A (portfolio date) stays the same for both rows.

# New row 2 (was row 3 / policy 2)
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 22890
$ws.Range("D2").Value = 42278
$ws.Range("E2").Value = 100000
$ws.Range("J2").Value = "N"
$ws.Range("K2").Value = 0.03
$ws.Range("L2").Value = 42278

# New row 3 (was row 2 / policy 1)
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 27873
$ws.Range("D3").Value = 44562
$ws.Range("E3").Value = 120000
$ws.Range("J3").Value = "S"
$ws.Range("K3").Value = 0.04
$ws.Range("L3").Value = 44562

# Row 4 was a leftover near-empty row; remove it entirely.
$ws.Rows.Item(4).Delete()

# Selection now covers the data rows, matching the saved view state.
$ws.Range("A2:L3").Select()
